$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# This document renders inline XML-ish markup (e.g. "<x>", "<exp>") as plain
# text runs styled to look like markup (blue/grey Courier New). The edit
# removes an outer "<x>...</x>" wrapper (an "element rendition spec") around
# two "<exp>ent</exp>" expansions, while keeping the inner "<exp>"/"</exp>"
# markup runs untouched. In one spot the word broken across the old "<x>"
# boundary ("ayse" + "m" + "<x>") needs to keep its trailing "m" as plain
# text instead of as part of the deleted markup run.
# ---------------------------------------------------------------------------

function Remove-LiteralSpan($doc, $searchText, $relStart, $relEnd) {
    # Finds the first literal occurrence of $searchText in the document and
    # deletes the sub-range [relStart, relEnd) (offsets relative to the
    # start of the match).
    $rng = $doc.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND:" $searchText
        return
    }
    $matchStart = $rng.Start
    $delRange = $doc.Range($matchStart + $relStart, $matchStart + $relEnd)
    Write-Host "Removing span:" $delRange.Text
    $delRange.Delete()
}

# 1) "...sable communem<x><exp>ent</exp></x> s'attaquera..."
#    -> remove the opening "<x>" right after "communem"
Remove-LiteralSpan $d "communem<x><exp>" 8 11

# 2) -> remove the closing "</x>" right after "</exp>", before " s'attaquera"
Remove-LiteralSpan $d "</exp></x> s'attaquera" 6 10

# 3) "...elle ne prend pas aysem<x><exp>ent</exp></x> sur les choses..."
#    -> the trailing "m" of "aysem" must survive as its own plain-text run
#       (same formatting as the surrounding prose) once "<x>" is removed.
$rng = $d.Content
$found = $rng.Find.Execute("aysem<x>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $matchEnd = $rng.End

    # Delete the "<x>" run (last 3 characters of the match).
    $xRange = $d.Range($matchEnd - 3, $matchEnd)
    Write-Host "Removing span:" $xRange.Text
    $xRange.Delete()

    # Split the trailing "m" off of "aysem" into its own run. A plain
    # property round-trip (toggle a trait on, then back off) forces the
    # engine to materialize "m" as a distinct run instead of leaving it
    # fused to "ayse", while its final formatting ends up identical to the
    # surrounding plain-text runs (color 000000).
    $mRange = $d.Range($matchEnd - 4, $matchEnd - 3)
    Write-Host "Splitting off trailing char:" $mRange.Text
    $mRange.Bold = 1
    $mRange = $d.Range($matchEnd - 4, $matchEnd - 3)
    $mRange.Bold = 0
} else {
    Write-Host "NOT FOUND: aysem<x>"
}

# 4) -> remove the closing "</x>" right after "</exp>", before " sur les choses"
Remove-LiteralSpan $d "</exp></x> sur les choses" 6 10
